$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 66; $r++) {
    $q = $ws.Cells.Item($r, 17).Value()   # Column Q = season_ending_year_y
    $s = $ws.Cells.Item($r, 19).Value()   # Column S = birth_year

    $ws.Cells.Item($r, 17).Value = $q - 1
    $ws.Cells.Item($r, 19).Value = $s + 1
}
